# Updates hashcode values in column B for specific rows, as part of an
# automated "hashcode" metadata refresh (data/metadata/hashcode.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    89  = "3a425473b901d99eeb2f8f05d1a7a9da"
    99  = "7332e19db9d80de1248db805e60f9312"
    110 = "a0cab0e46f110ea81f706b2fc5953f20"
    154 = "dadb7be999dbd63f806299bfafbc6261"
    160 = "25264021f32130c246ff1dcdeec483d0"
    281 = "181895aa68478a8ce5e37e3a6123fdf6"
    338 = "c16252edd9bbad81bece7e1e437aeca5"
    511 = "b3c0471f6ab03fe79ed3515cd46b22cc"
    545 = "6872b106d46507f66af37d33523f76f9"
    559 = "a43aad2a42277be6fc85233bafe81f21"
    565 = "2ba2af195a7150411e9edbf214040e44"
    596 = "db79560a07b943a028661bf9ac58f8cf"
    677 = "16b63d480f3d50d78a869c19ab998727"
    712 = "32cabfb6d54c47197f02bfa132f2bceb"
    780 = "7b32c2e2138ad20d6de90800ca768f42"
    823 = "1240d1925d5bb6781d888325f1408e49"
    827 = "18959c8132fbe58132b63e2ed262ede7"
    828 = "683ad9d5a62eedccab952d06bed5a4f7"
    837 = "c23d1d2e9e89bd032e026d27dfcc8827"
    839 = "97010d418992034607b9ffb8ac4a8020"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
